# corrected merge and export to excel WITH index
#
# Two changes to the df_merged sheet's data rows (2-12):
#   1. Column A ("year") is re-exported as text instead of a number, using
#      the same (bordered / bold / centered) style already used by the
#      header row.
#   2. Every cell that was a bare 0 placeholder (the pandas NaN-as-0 export
#      artifact) is cleared to an empty string instead of a literal zero,
#      while keeping the cell's original (default) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$years = @{
    2  = "2021"
    3  = "2020"
    4  = "2019"
    5  = "2018"
    6  = "2017"
    7  = "2016"
    8  = "2015"
    9  = "2014"
    10 = "2013"
    11 = "2012"
    12 = "2011"
}

foreach ($r in $years.Keys) {
    # Leading apostrophe forces the numeric-looking year to be stored as text.
    $ws.Cells.Item($r, 1).Value = "'" + $years[$r]
}

# Re-apply the header's formatting (bold / bordered / centered, style index
# used by row 1) onto the year column so it keeps looking like the header
# once it becomes text.
$ws.Range("A1").Copy()
$ws.Range("A2:A12").PasteSpecial(-4122)

$zeroCells = @(
    "E2","F2","G2","H2","M2",
    "E3","F3","G3","I3","J3","K3","M3",
    "E4","F4","G4","H4","I4","J4","M4",
    "E5","F5","G5","H5","J5","K5",
    "E6","F6","G6","H6","J6","L6",
    "E7","F7","G7","I7","J7","M7",
    "E8","F8","H8","I8","J8","L8","M8",
    "E9","F9","G9","H9","K9","L9",
    "E10","F10","G10","H10","J10","M10",
    "E11","F11","G11","H11","I11","J11","K11","L11","M11",
    "E12","F12","G12","H12","I12","K12","M12"
)

foreach ($addr in $zeroCells) {
    # Leading apostrophe forces an explicit (empty) text cell rather than
    # clearing the cell outright.
    $ws.Range($addr).Value = "'"
}

# Restore each cleared cell's original (unstyled) formatting - setting a
# text value above can otherwise stamp a quote-prefixed style onto it.
$ws.Range("B2").Copy()
foreach ($addr in $zeroCells) {
    $ws.Range($addr).PasteSpecial(-4122)
}
